# "carga trucha guido mirala" - marcar como completados (valor 1) los
# parametros D que faltaban para Banjo (fila 7) y Organ (fila 10).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("D7").Value = 1
$ws.Range("D10").Value = 1

# Dejar el cursor/seleccion donde quedo Guido tras cargar los datos
$ws.Range("J36:J37").Select()
